$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16 and 17 (roster entries #14 and #15) trade places: "Julian
# Champagnie (TW)" and "Sandro Mamukelashvili" swap everything (jersey #,
# position, height, weight, birth date, nationality, experience, college,
# bbref url) except column A (the "No." index), which stays sequential:
# 14 stays on row 16, 15 stays on row 17.
#
# Use Copy/Paste (via a scratch row far below the data) so each cell keeps
# its original type/format (e.g. the text "1" in the Exp column doesn't get
# reinterpreted as a number) instead of manually re-typing values.

$scratch = $ws.Range("B200:K200")

$ws.Range("B16:K16").Copy($scratch)
$ws.Range("B17:K17").Copy($ws.Range("B16:K16"))
$scratch.Copy($ws.Range("B17:K17"))
$scratch.Clear()
